$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply an integer number format to the whole column D (this is what produced
# the new combined cell-format entries in styles.xml: one for the already
# bordered header cell D7, one as the plain column-level style).
$colD = $ws.Range("D1:D1048576")
$colD.NumberFormat = "0"
$colD.Select()

# Clear the header text that was previously in row 7 (the shared strings -
# "#", "Numero Documento", ... - are gone from the saved workbook), while
# keeping the cell formatting/borders in place.
$ws.Range("A7:I7").ClearContents()
